$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update formulas/values on row 2
$ws.Range("B2").Formula = "=SUM(E5:E20)"
$ws.Range("C2").Formula = "=SUM(D5:D20)"
$ws.Range("D2").ClearContents()

# Update the selected range / active cell
$ws.Range("F3").Select()
